$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "69.584.26"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -1.62%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.483.84"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -4.08%  "

$ws.Range("E4").Value = "  +0.07%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "579.58"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -4.25%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "193.10"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -2.89%  "

$ws.Range("E7").Value = "  -2.45%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "3.472.30"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -3.98%  "

$ws.Range("E9").Value = "  -0.03%  "

$ws.Range("E10").Value = "  -7.75%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.618"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -4.42%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "51.47"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -4.39%  "

$ws.Range("E13").Value = "  -6.40%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "9.15"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -4.18%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "4.034.33"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -4.09%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "649.71"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -3.69%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "69.449.45"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -2.01%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "3.478.68"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -4.34%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "12.32"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -5.38%  "

$ws.Range("E20").Value = "  -1.74%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "18.23"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -4.42%  "

$ws.Range("E22").Value = "  -5.21%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "18.17"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -2.26%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "5.27"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -1.89%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "98.88"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -6.30%  "

$ws.Range("E26").Value = "  -7.39%  "

$ws.Range("E27").Value = "  -3.85%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "10.03"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -4.02%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "9.35"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -4.74%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "32.59"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -4.23%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "4.28"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -8.75%  "

$ws.Range("E32").Value = "  -6.34%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "11.61"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -4.85%  "

$ws.Range("E34").Value = "  -5.09%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "60.95"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -3.90%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "3.724.78"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -6.14%  "

$ws.Range("E37").Value = "  +0.14%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "521.13"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +3.00%  "

$ws.Range("E39").Value = "  -9.03%  "

$ws.Range("E40").Value = "  -3.47%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "3.51"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.83%  "

$ws.Range("E42").Value = "  -3.72%  "

$ws.Range("B43").Value = "CoreDAO"
$ws.Range("C43").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "3.55"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +71.51%  "

$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.133"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -2.55%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "34.35"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -6.70%  "

$ws.Range("E46").Value = "  -4.05%  "

$ws.Range("E47").Value = "  -4.36%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.84"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -7.41%  "

$ws.Range("E49").Value = "  -4.27%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.42%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "8.17"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -5.83%  "
